# Scheduled-runner update of Masamune_Profits price/profit figures.
# Refreshes currentAveragePrice(NQ/HQ), LevePrice(NQ/HQ) and LeveProfit(NQ/HQ)
# columns (H:N) on several leve rows across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR
# sheets to reflect newly-pulled market data.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H13").Value = 10002.5
$ws.Range("I13").Value = 9999
$ws.Range("J13").Value = 10006
$ws.Range("K13").Value = 9999
$ws.Range("L13").Value = 10006
$ws.Range("M13").Value = -9830
$ws.Range("N13").Value = -10344

$ws.Range("H28").Value = 466
$ws.Range("I28").Value = 256.36365
$ws.Range("J28").Value = 927.2
$ws.Range("K28").Value = 256.36365
$ws.Range("L28").Value = 927.2
$ws.Range("M28").Value = 228.63635
$ws.Range("N28").Value = -1897.2

$ws.Range("H76").Value = 3047.4036
$ws.Range("I76").Value = 3030.9456
$ws.Range("K76").Value = 3030.9456
$ws.Range("M76").Value = -2715.9456

$ws.Range("H79").Value = 3047.4036
$ws.Range("I79").Value = 3030.9456
$ws.Range("K79").Value = 3030.9456
$ws.Range("M79").Value = -1938.9456

$ws.Range("H98").Value = 38337.145
$ws.Range("I98").Value = 1135.3572
$ws.Range("J98").Value = 112740.71
$ws.Range("K98").Value = 1135.3572
$ws.Range("L98").Value = 112740.71
$ws.Range("M98").Value = 362.6428000000001
$ws.Range("N98").Value = -115736.71

$ws.Range("H112").Value = 1297.75
$ws.Range("J112").Value = 1297.5454
$ws.Range("L112").Value = 3892.6362
$ws.Range("N112").Value = -6108.6362

$ws.Range("H122").Value = 38337.145
$ws.Range("I122").Value = 1135.3572
$ws.Range("J122").Value = 112740.71
$ws.Range("K122").Value = 3406.0716
$ws.Range("L122").Value = 338222.13
$ws.Range("M122").Value = -956.0715999999998
$ws.Range("N122").Value = -343122.13

$ws.Range("H141").Value = 3162.3572
$ws.Range("I141").Value = 3034.4
$ws.Range("J141").Value = 3233.4443
$ws.Range("K141").Value = 9103.2
$ws.Range("L141").Value = 9700.332900000001
$ws.Range("M141").Value = -3923.200000000001
$ws.Range("N141").Value = -20060.3329

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2452.2917
$ws.Range("I61").Value = 1964.1177
$ws.Range("J61").Value = 3637.8572
$ws.Range("K61").Value = 1964.1177
$ws.Range("L61").Value = 3637.8572
$ws.Range("M61").Value = -1752.1177
$ws.Range("N61").Value = -4061.8572

$ws.Range("H76").Value = 25000
$ws.Range("J76").Value = 25000
$ws.Range("L76").Value = 25000
$ws.Range("N76").Value = -25676

$ws.Range("H79").Value = 25000
$ws.Range("J79").Value = 25000
$ws.Range("L79").Value = 25000
$ws.Range("N79").Value = -27340

$ws.Range("H94").Value = 30000
$ws.Range("J94").Value = 30000
$ws.Range("L94").Value = 30000
$ws.Range("N94").Value = -31802

$ws.Range("H134").Value = 48000
$ws.Range("J134").Value = 48000
$ws.Range("L134").Value = 48000
$ws.Range("N134").Value = -58140

$ws.Range("H136").Value = 2452.2917
$ws.Range("I136").Value = 1964.1177
$ws.Range("J136").Value = 3637.8572
$ws.Range("K136").Value = 5892.3531
$ws.Range("L136").Value = 10913.5716
$ws.Range("M136").Value = -3342.3531
$ws.Range("N136").Value = -16013.5716

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 502500000
$ws.Range("I5").Value = 502500000
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 502500000
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -502499887
$ws.Range("N5").ClearContents()

$ws.Range("H107").Value = 2140.5715
$ws.Range("I107").Value = 2159.5454
$ws.Range("J107").Value = 2071
$ws.Range("K107").Value = 2159.5454
$ws.Range("L107").Value = 2071
$ws.Range("M107").Value = -239.5454
$ws.Range("N107").Value = -5911

$ws.Range("H134").Value = 1771.8125
$ws.Range("I134").Value = 1204.25
$ws.Range("J134").Value = 3474.5
$ws.Range("K134").Value = 3612.75
$ws.Range("L134").Value = 10423.5
$ws.Range("M134").Value = -1077.75
$ws.Range("N134").Value = -15493.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 12501440
$ws.Range("I58").Value = 803.0769
$ws.Range("J58").Value = 35716908
$ws.Range("K58").Value = 803.0769
$ws.Range("L58").Value = 35716908
$ws.Range("M58").Value = -600.0769
$ws.Range("N58").Value = -35717314

$ws.Range("H99").Value = 2745
$ws.Range("I99").Value = 2696.3635
$ws.Range("K99").Value = 2696.3635
$ws.Range("M99").Value = -1198.3635

$ws.Range("H126").Value = 2745
$ws.Range("I126").Value = 2696.3635
$ws.Range("K126").Value = 8089.0905
$ws.Range("M126").Value = -5619.0905

$ws.Range("H136").Value = 12501440
$ws.Range("I136").Value = 803.0769
$ws.Range("J136").Value = 35716908
$ws.Range("K136").Value = 2409.2307
$ws.Range("L136").Value = 107150724
$ws.Range("M136").Value = 140.7692999999999
$ws.Range("N136").Value = -107155824

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 2942653.2
$ws.Range("I92").Value = 11765462
$ws.Range("J92").Value = 1716.8667
$ws.Range("K92").Value = 35296386
$ws.Range("L92").Value = 5150.6001
$ws.Range("M92").Value = -35295138
$ws.Range("N92").Value = -7646.6001

$ws.Range("H131").Value = 901.88
$ws.Range("I131").Value = 596.5
$ws.Range("J131").Value = 928.43475
$ws.Range("K131").Value = 1789.5
$ws.Range("L131").Value = 2785.30425
$ws.Range("M131").Value = 3250.5
$ws.Range("N131").Value = -12865.30425

$ws.Range("H137").Value = 12019.571
$ws.Range("I137").Value = 3441
$ws.Range("J137").Value = 19818.273
$ws.Range("K137").Value = 10323
$ws.Range("L137").Value = 59454.819
$ws.Range("M137").Value = -5223
$ws.Range("N137").Value = -69654.819

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H74").Value = 40771.43
$ws.Range("J74").Value = 40771.43
$ws.Range("L74").Value = 40771.43
$ws.Range("N74").Value = -42643.43

$ws.Range("H77").Value = 40771.43
$ws.Range("J77").Value = 40771.43
$ws.Range("L77").Value = 122314.29
$ws.Range("N77").Value = -131674.29

$ws.Range("H102").Value = 1321.9062
$ws.Range("I102").Value = 930.04
$ws.Range("J102").Value = 2721.4285
$ws.Range("K102").Value = 930.04
$ws.Range("L102").Value = 2721.4285
$ws.Range("M102").Value = 691.96
$ws.Range("N102").Value = -5965.4285

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("M76").ClearContents()

$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("M79").ClearContents()

$ws.Range("H82").Value = 5956041
$ws.Range("I82").Value = 3436.4
$ws.Range("J82").Value = 9263044
$ws.Range("K82").Value = 3436.4
$ws.Range("L82").Value = 9263044
$ws.Range("M82").Value = -3075.4
$ws.Range("N82").Value = -9263766

$ws.Range("H85").Value = 5956041
$ws.Range("I85").Value = 3436.4
$ws.Range("J85").Value = 9263044
$ws.Range("K85").Value = 3436.4
$ws.Range("L85").Value = 9263044
$ws.Range("M85").Value = -2188.4
$ws.Range("N85").Value = -9265540

$ws.Range("H93").Value = 1133.2
$ws.Range("I93").Value = 500
$ws.Range("J93").Value = 1178.4286
$ws.Range("K93").Value = 500
$ws.Range("L93").Value = 1178.4286
$ws.Range("M93").Value = 748
$ws.Range("N93").Value = -3674.4286

$ws.Range("H132").Value = 3879.36
$ws.Range("I132").Value = 2656.8572
$ws.Range("J132").Value = 5435.273
$ws.Range("K132").Value = 7970.571599999999
$ws.Range("L132").Value = 16305.819
$ws.Range("M132").Value = -5440.571599999999
$ws.Range("N132").Value = -21365.819

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 12822011
$ws.Range("I132").Value = 1333.2
$ws.Range("J132").Value = 26317462
$ws.Range("K132").Value = 3999.6
$ws.Range("L132").Value = 78952386
$ws.Range("M132").Value = -1469.6
$ws.Range("N132").Value = -78957446
